$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 3.272327238179451
$ws.Range("C2").Value = 1.626987699542094
$ws.Range("D2").Value = 3.223369029078222
$ws.Range("E2").Value = 0.5333859586016987
$ws.Range("G2").Value = 8.656069925401464

# Row 3
$ws.Range("B3").Value = 3.272327238179451
$ws.Range("C3").Value = 1.626987699542094
$ws.Range("D3").Value = 3.223369029078222
$ws.Range("E3").Value = 13.86384647080068
$ws.Range("G3").Value = 21.98653043760045

# Row 4
$ws.Range("B4").Value = 0.04172184405617529
$ws.Range("C4").Value = 9.983522426115931
$ws.Range("D4").Value = 0.1496068669990043
$ws.Range("E4").Value = 13.86384647080068
$ws.Range("G4").Value = 24.03869760797179
